$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the task name in A2 ("ahahah123123123" -> "emre bey")
$ws.Range("A2").Value = "emre bey"

# Remove the now-obsolete row 3 ("selam" / "Continue...")
$ws.Rows(3).Delete()
